{"js": "const table = context.document.body.tables.getFirst();\nconst replacements = [\n  { row: 0, col: 0, oldText: \"44\u00f76=7, 2\", newText: \"37\u00f72=18, 1\" },\n  { row: 0, col: 1, oldText: \"48\u00f79=5, 3\", newText: \"98\u00f74=24, 2\" },\n  { row: 0, col: 2, oldText: \"78\u00f79=8, 6\", newText: \"40\u00f73=13, 1\" },\n  { row: 0, col: 3, oldText: \"30\u00f76=5, 0\", newText: \"37\u00f79=4, 1\" },\n  { row: 0, col: 4, oldText: \"83\u00f75=16, 3\", newText: \"41\u00f77=5, 6\" },\n  { row: 4, col: 0, oldText: \"38\u00f73=12, 2\", newText: \"11\u00f79=1, 2\" },\n  { row: 4, col: 1, oldText: \"97\u00f79=10, 7\", newText: \"92\u00f74=23, 0\" },\n  { row: 4, col: 2, oldText: \"63\u00f75=12, 3\", newText: \"23\u00f72=11, 1\" },\n  { row: 4, col: 3, oldText: \"33\u00f74=8, 1\", newText: \"39\u00f74=9, 3\" },\n  { row: 4, col: 4, oldText: \"54\u00f78=6, 6\", newText: \"30\u00f79=3, 3\" },\n  { row: 8, col: 0, oldText: \"54\u00f75=10, 4\", newText: \"28\u00f72=14, 0\" },\n  { row: 8, col: 1, oldText: \"33\u00f72=16, 1\", newText: \"80\u00f72=40, 0\" },\n  { row: 8, col: 2, oldText: \"29\u00f78=3, 5\", newText: \"79\u00f79=8, 7\" },\n  { row: 8, col: 3, oldText: \"34\u00f74=8, 2\", newText: \"21\u00f73=7, 0\" },\n  { row: 8, col: 4, oldText: \"88\u00f76=14, 4\", newText: \"43\u00f72=21, 1\" },\n  { row: 12, col: 0, oldText: \"96\u00f74=24, 0\", newText: \"15\u00f76=2, 3\" },\n  { row: 12, col: 1, oldText: \"96\u00f78=12, 0\", newText: \"53\u00f77=7, 4\" },\n  { row: 12, col: 2, oldText: \"74\u00f74=18, 2\", newText: \"45\u00f73=15, 0\" },\n  { row: 12, col: 3, oldText: \"83\u00f79=9, 2\", newText: \"97\u00f73=32, 1\" },\n  { row: 12, col: 4, oldText: \"46\u00f76=7, 4\", newText: \"53\u00f74=13, 1\" },\n  { row: 16, col: 0, oldText: \"22\u00f73=7, 1\", newText: \"63\u00f72=31, 1\" },\n  { row: 16, col: 1, oldText: \"26\u00f75=5, 1\", newText: \"25\u00f77=3, 4\" },\n  { row: 16, col: 2, oldText: \"94\u00f73=31, 1\", newText: \"35\u00f78=4, 3\" },\n  { row: 16, col: 3, oldText: \"37\u00f72=18, 1\", newText: \"81\u00f78=10, 1\" },\n  { row: 16, col: 4, oldText: \"22\u00f78=2, 6\", newText: \"31\u00f76=5, 1\" },\n];\n\n// Locate the specific run of text inside each target cell and replace it,\n// preserving paragraph/run formatting (search+insertText \"Replace\").\nconst foundRanges = [];\nfor (const r of replacements) {\n  const cell = table.getCell(r.row, r.col);\n  const results = cell.body.search(r.oldText, { matchCase: true });\n  results.load(\"items\");\n  foundRanges.push({ results, newText: r.newText });\n}\nawait context.sync();\n\nfor (const { results, newText } of foundRanges) {\n  results.items[0].insertText(newText, \"Replace\");\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# (row, col, expectedOldText, newText) - 1-indexed table cell coordinates\n$replacements = @(\n    @(1, 1, \"44\u00f76=7, 2\", \"37\u00f72=18, 1\"),\n    @(1, 2, \"48\u00f79=5, 3\", \"98\u00f74=24, 2\"),\n    @(1, 3, \"78\u00f79=8, 6\", \"40\u00f73=13, 1\"),\n    @(1, 4, \"30\u00f76=5, 0\", \"37\u00f79=4, 1\"),\n    @(1, 5, \"83\u00f75=16, 3\", \"41\u00f77=5, 6\"),\n    @(5, 1, \"38\u00f73=12, 2\", \"11\u00f79=1, 2\"),\n    @(5, 2, \"97\u00f79=10, 7\", \"92\u00f74=23, 0\"),\n    @(5, 3, \"63\u00f75=12, 3\", \"23\u00f72=11, 1\"),\n    @(5, 4, \"33\u00f74=8, 1\", \"39\u00f74=9, 3\"),\n    @(5, 5, \"54\u00f78=6, 6\", \"30\u00f79=3, 3\"),\n    @(9, 1, \"54\u00f75=10, 4\", \"28\u00f72=14, 0\"),\n    @(9, 2, \"33\u00f72=16, 1\", \"80\u00f72=40, 0\"),\n    @(9, 3, \"29\u00f78=3, 5\", \"79\u00f79=8, 7\"),\n    @(9, 4, \"34\u00f74=8, 2\", \"21\u00f73=7, 0\"),\n    @(9, 5, \"88\u00f76=14, 4\", \"43\u00f72=21, 1\"),\n    @(13, 1, \"96\u00f74=24, 0\", \"15\u00f76=2, 3\"),\n    @(13, 2, \"96\u00f78=12, 0\", \"53\u00f77=7, 4\"),\n    @(13, 3, \"74\u00f74=18, 2\", \"45\u00f73=15, 0\"),\n    @(13, 4, \"83\u00f79=9, 2\", \"97\u00f73=32, 1\"),\n    @(13, 5, \"46\u00f76=7, 4\", \"53\u00f74=13, 1\"),\n    @(17, 1, \"22\u00f73=7, 1\", \"63\u00f72=31, 1\"),\n    @(17, 2, \"26\u00f75=5, 1\", \"25\u00f77=3, 4\"),\n    @(17, 3, \"94\u00f73=31, 1\", \"35\u00f78=4, 3\"),\n    @(17, 4, \"37\u00f72=18, 1\", \"81\u00f78=10, 1\"),\n    @(17, 5, \"22\u00f78=2, 6\", \"31\u00f76=5, 1\"),\n)\n\nforeach ($rep in $replacements) {\n    $row = $rep[0]; $col = $rep[1]; $oldText = $rep[2]; $newText = $rep[3]\n    $cell = $t.Cell($row, $col)\n    # Cell.Range.Text includes the trailing cell-mark (CR + BEL); trim before comparing.\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $oldText) {\n        Write-Output \"warning: cell ($row,$col) expected '$oldText' but found '$current'\"\n    }\n    $cell.Range.Text = $newText\n}"}
